# Refresh the crypto price / 1h-volume table (columns D and E) with the
# latest scraped values. Every cell on this sheet is stored as text, so a
# leading apostrophe is used on the "Price" column where the new value
# would otherwise look like a plain number to Excel's input parser (which
# would silently convert it to a numeric cell instead of text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.994.80"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "3.121.65"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'586.51"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "'146.51"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.115.58"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  +9.39%  "
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").Value = "'0.464"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").Value = "'37.34"
$ws.Range("E14").Value = "  +4.51%  "
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").Value = "3.638.15"
$ws.Range("D17").Value = "63.860.70"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").Value = "3.120.12"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "'464.70"
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "'7.54"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").Value = "'13.17"
$ws.Range("E24").Value = "  -3.43%  "
$ws.Range("D25").Value = "'81.72"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'8.96"
$ws.Range("E27").Value = "  +8.83%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "'6.89"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").Value = "'26.98"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("E34").Value = "  +6.47%  "
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("E37").Value = "  +9.52%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").Value = "'447.52"
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "'0.0372"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "2.880.43"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").Value = "'35.73"
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("D49").Value = "'123.38"
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "'24.66"
$ws.Range("E51").Value = "  -1.27%  "
